# Update "F" column ("想去人数" / "want-to-go count") figures that changed
# between the previous and newly generated gh-pages data snapshot.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 3420
$ws.Range("F4").Value  = 586
$ws.Range("F5").Value  = 840
$ws.Range("F6").Value  = 320
$ws.Range("F7").Value  = 275
$ws.Range("F9").Value  = 165
$ws.Range("F10").Value = 642
$ws.Range("F11").Value = 197
$ws.Range("F12").Value = 449
$ws.Range("F13").Value = 74
$ws.Range("F14").Value = 495
$ws.Range("F15").Value = 349
$ws.Range("F16").Value = 61
$ws.Range("F18").Value = 100
$ws.Range("F19").Value = 188

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 7
$ws.Range("F9").Value = 180

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6211
$ws.Range("F4").Value = 747
$ws.Range("F5").Value = 1794
$ws.Range("F6").Value = 118

# Sheet "全部类型" (All types) - union of the three sheets above
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 6211
$ws.Range("F4").Value  = 747
$ws.Range("F5").Value  = 1794
$ws.Range("F6").Value  = 3420
$ws.Range("F7").Value  = 118
$ws.Range("F9").Value  = 586
$ws.Range("F10").Value = 840
$ws.Range("F11").Value = 320
$ws.Range("F12").Value = 275
$ws.Range("F17").Value = 165
$ws.Range("F19").Value = 7
$ws.Range("F20").Value = 642
$ws.Range("F22").Value = 197
$ws.Range("F24").Value = 449
$ws.Range("F25").Value = 180
$ws.Range("F26").Value = 74
$ws.Range("F27").Value = 495
$ws.Range("F29").Value = 349
$ws.Range("F30").Value = 61
$ws.Range("F34").Value = 100
$ws.Range("F40").Value = 188
